$wb = $excel.ActiveWorkbook

# --- Overview sheet: refresh "Latest HO Xliff Generate Date" for the
#     19f5f858... report row (rows 4 & 5 both point at that file) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-31 10:18:09"
$wsOverview.Range("G5").Value = "2016-08-31 10:18:09"

# --- zh-cn sheet: Priority changed from "ht" to "mt" for the
#     19f5f858... file, its Correspond Handoff Datetime and both
#     rows' Correspond Handback DateTime refreshed ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-31 10:17:57"
$wsZhCn.Range("H5").Value = "2016-08-31 10:17:57"
$wsZhCn.Range("K4").Value = "2016-08-31 10:18:31"
$wsZhCn.Range("K5").Value = "2016-08-31 10:18:31"

# --- de-de sheet: Priority changed from "ht" to "mt" for the
#     19f5f858... file, and its Correspond Handback DateTime refreshed ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("K4").Value = "2016-08-31 10:18:38"
$wsDeDe.Range("K5").Value = "2016-08-31 10:18:38"
